$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Update the "distance" sheet: insert the new "trajectory-orange-2" row
#    and column into the distance matrix, expanding it from a 5x5 table
#    (A1:F6) to a 6x6 table (A1:G7).
# ---------------------------------------------------------------------------
$wsDistance = $wb.Worksheets.Item("distance")

# Insert a new row above row 4 (the old row for "trajectory-red-1"), and a
# new column before column D (the old column for "trajectory-red-1"), to make
# room for the "trajectory-orange-2" row/column in sorted order.
$wsDistance.Range("A4").EntireRow.Insert()
$wsDistance.Range("D1").EntireColumn.Insert()

# Header row (row 1) and row labels (column A)
$wsDistance.Range("D1").Value = "trajectory-orange-2"
$wsDistance.Range("A4").Value = "trajectory-orange-2"

# Full (symmetric) distance matrix, rows/cols 2-7 = B:G correspond to
# trajectory-green-1, trajectory-orange-1, trajectory-orange-2,
# trajectory-red-1, trajectory-red-2, trajectory-red-3
$wsDistance.Range("B2").Value = 0
$wsDistance.Range("C2").Value = 1.83458544342631
$wsDistance.Range("D2").Value = 1.81558152790223
$wsDistance.Range("E2").Value = 0.830278087225997
$wsDistance.Range("F2").Value = 1.00815783653297
$wsDistance.Range("G2").Value = 0.943492689191259

$wsDistance.Range("B3").Value = 1.83458544342631
$wsDistance.Range("C3").Value = 0
$wsDistance.Range("D3").Value = 0.130761554072711
$wsDistance.Range("E3").Value = 1.25554951954802
$wsDistance.Range("F3").Value = 1.1394684404619
$wsDistance.Range("G3").Value = 1.19512366294384

$wsDistance.Range("B4").Value = 1.81558152790223
$wsDistance.Range("C4").Value = 0.130761554072711
$wsDistance.Range("D4").Value = 0
$wsDistance.Range("E4").Value = 1.22163041053492
$wsDistance.Range("F4").Value = 1.10097133040577
$wsDistance.Range("G4").Value = 1.15938604841685

$wsDistance.Range("B5").Value = 0.830278087225997
$wsDistance.Range("C5").Value = 1.25554951954802
$wsDistance.Range("D5").Value = 1.22163041053492
$wsDistance.Range("E5").Value = 0
$wsDistance.Range("F5").Value = 0.178268246647196
$wsDistance.Range("G5").Value = 0.103527405716764

$wsDistance.Range("B6").Value = 1.00815783653297
$wsDistance.Range("C6").Value = 1.1394684404619
$wsDistance.Range("D6").Value = 1.10097133040577
$wsDistance.Range("E6").Value = 0.178268246647196
$wsDistance.Range("F6").Value = 0
$wsDistance.Range("G6").Value = 0.115738277320899

$wsDistance.Range("B7").Value = 0.943492689191259
$wsDistance.Range("C7").Value = 1.19512366294384
$wsDistance.Range("D7").Value = 1.15938604841685
$wsDistance.Range("E7").Value = 0.103527405716764
$wsDistance.Range("F7").Value = 0.115738277320899
$wsDistance.Range("G7").Value = 0

# Column widths for the new A1:G7 table (closest character-width input that
# rounds, under this engine's pixel-based column-width model, to the
# target stored widths of 18.05 and 18.9)
$wsDistance.Range("A:A").ColumnWidth = 17.1666666667
$wsDistance.Range("B:B").ColumnWidth = 18
$wsDistance.Range("C:C").ColumnWidth = 17.1666666667
$wsDistance.Range("D:G").ColumnWidth = 18

# ---------------------------------------------------------------------------
# 2. Restore each sheet's selection to a single cell/range (cleaning up the
#    leftover "B4:G4" selection remnants) and make "distance" the active
#    sheet (it was "sorted" before).
# ---------------------------------------------------------------------------
$wsTrajectory = $wb.Worksheets.Item("trajectory")
$wsTrajectory.Activate()
$wsTrajectory.Range("A8").Select()

$wsGenotype = $wb.Worksheets.Item("genotype")
$wsGenotype.Activate()
$wsGenotype.Range("B3").Select()

$wsSorted = $wb.Worksheets.Item("sorted")
$wsSorted.Activate()
$wsSorted.Range("B4").Select()

$wsEdges = $wb.Worksheets.Item("edges")
$wsEdges.Activate()
$wsEdges.Range("A5").Select()

$wsDistance.Activate()
$wsDistance.Range("A1").Select()
